$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("list")

# --- Update the enum value / description pairs for order_status_type,
#     order_priority_type and order_category_type (rows 3-5), and refresh
#     the asset_category_type pair on row 1 to its new shared-string slot.
$ws.Range("B1").Value = "('F', 'E', 'A')"
$ws.Range("C1").Value = "facility, equipment, appliances"

$ws.Range("B3").Value = "('CAN', 'NEG', 'PEN', 'SUS', 'FIL', 'EXE', 'CON')"
$ws.Range("C3").Value = "cancelada, negada, pendente, suspenso, fila de espera, execução, concluída"

$ws.Range("B4").Value = "('BAI', 'NOR', 'ALT', 'URG')"
$ws.Range("C4").Value = "baixa, normal, alta, urgente"

$ws.Range("B5").Value = "('EST', 'FOR', 'INF', 'ELE', 'HID', 'MAR', 'PIS', 'REV', 'VED', 'VID', 'SER')"
$ws.Range("C5").Value = "avaliação estrutural, reparo em forro, infiltração, instalações elétricas, instalações hidrossanitárias, marcenaria, reparo em piso, revestimento, vedação espacial, vidraçaria/esquadria, serralheria"

# --- Wrap the long description column so the text is readable, which
#     introduces the new cellXfs style used by column C.
$ws.Range("C1:C5").WrapText = $true

# --- Widen columns B (values) and C (descriptions) to fit their new,
#     longer content.
$ws.Columns.Item(2).ColumnWidth = 60.42578125
$ws.Columns.Item(3).ColumnWidth = 69.7109375

# --- Row 5 now needs extra height to show the long, wrapped category list.
$ws.Rows.Item(5).RowHeight = 45

# --- Leave the selection where the author left it when they finished.
$ws.Range("C8").Select()
